$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$helper = $ws.Range("Z1")

$ws.Range("D2").Value = "59.986.77"
$ws.Range("E2").Value = "  +2.23%  "
$ws.Range("D3").Value = "3.189.55"
$ws.Range("E3").Value = "  +0.93%  "
$helper.Formula = '="536.67"'
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +1.28%  "
$helper.Formula = '="144.95"'
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +3.59%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -2.48%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("E11").Value = "  -1.26%  "
$ws.Range("D12").Value = "3.739.74"
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("E13").Value = "  -2.78%  "
$helper.Formula = '="25.72"'
$helper.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").Value = "59.985.51"
$ws.Range("E16").Value = "  +2.14%  "
$ws.Range("D17").Value = "3.172.25"
$ws.Range("E17").Value = "  +1.30%  "
$helper.Formula = '="6.23"'
$helper.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  -0.50%  "
$helper.Formula = '="13.22"'
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  +1.63%  "
$ws.Range("E20").Value = "  +0.75%  "
$helper.Formula = '="368.75"'
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("E23").Value = "  -1.92%  "
$helper.Formula = '="69.51"'
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("E25").Value = "  +1.64%  "
$helper.Formula = '="8.60"'
$helper.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  +3.47%  "
$helper.Formula = '="0.990"'
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -0.92%  "
$ws.Range("E28").Value = "  +0.51%  "
$helper.Formula = '="22.45"'
$helper.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("E30").Value = "  +0.22%  "
$helper.Formula = '="6.09"'
$helper.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  +0.46%  "
$helper.Formula = '="5.27"'
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  +2.82%  "
$ws.Range("E33").Value = "  +4.92%  "
$ws.Range("E34").Value = "  +2.55%  "
$helper.Formula = '="157.89"'
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("E36").Value = "  +1.82%  "
$helper.Formula = '="26.32"'
$helper.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +5.25%  "
$ws.Range("D38").Value = "2.785.66"
$ws.Range("E38").Value = "  +5.23%  "
$ws.Range("E39").Value = "  +2.78%  "
$helper.Formula = '="0.0309"'
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  +7.03%  "
$ws.Range("E41").Value = "  -0.26%  "
$helper.Formula = '="4.20"'
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -1.87%  "
$helper.Formula = '="39.94"'
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  +2.04%  "
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("D46").Value = "3.230.60"
$ws.Range("E46").Value = "  +0.94%  "
$helper.Formula = '="0.979"'
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("E48").Value = "  -0.94%  "
$helper.Formula = '="20.55"'
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("E50").Value = "  +5.64%  "

$helper.ClearContents()
$excel.CutCopyMode = $false
